$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old last data row (row 9) entirely so the table shrinks from 8 to 7 data rows
$ws.Range("A9:G9").Delete()

# Update header row (row 1) with new column headers, adding mean/std split columns
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "One Year Alt mean"
$ws.Range("D1").Value = "One Year Alt std"
$ws.Range("E1").Value = "Two Year Alt mean"
$ws.Range("F1").Value = "Two Year Alt std"
$ws.Range("G1").Value = "Three Year Alt mean"
$ws.Range("H1").Value = "Three Year Alt std"
$ws.Range("I1").Value = "Five Year Alt mean"
$ws.Range("J1").Value = "Five Year Alt std"
$ws.Range("K1").Value = "Ten Year Alt mean"
$ws.Range("L1").Value = "Ten Year Alt std"

# Copy the existing header formatting (bold font + border + center/top alignment)
# from G1 onto the newly added header cells H1:L1
$ws.Range("G1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)

# Row 2: LR
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.9080895633060564
$ws.Range("D2").Value = 0.009108492687851251
$ws.Range("E2").Value = 0.8966449202536892
$ws.Range("F2").Value = 0.008260723024058286
$ws.Range("G2").Value = 0.8874894370078259
$ws.Range("H2").Value = 0.01297015642625823
$ws.Range("I2").Value = 0.8781015495216818
$ws.Range("J2").Value = 0.0178074543301257
$ws.Range("K2").Value = 0.8703992760999887
$ws.Range("L2").Value = 0.01123659427447213

# Row 3: LDA
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.910854055771828
$ws.Range("D3").Value = 0.008813618922141053
$ws.Range("E3").Value = 0.9010305415027169
$ws.Range("F3").Value = 0.008235299189123073
$ws.Range("G3").Value = 0.8941401418196515
$ws.Range("H3").Value = 0.01129540148019561
$ws.Range("I3").Value = 0.8864302285616084
$ws.Range("J3").Value = 0.02053125033837994
$ws.Range("K3").Value = 0.8782462391132224
$ws.Range("L3").Value = 0.0141767843408919

# Row 4: KNN
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.8550544212282514
$ws.Range("D4").Value = 0.007959502690768923
$ws.Range("E4").Value = 0.8607326781019868
$ws.Range("F4").Value = 0.01213033894896635
$ws.Range("G4").Value = 0.8584369220971672
$ws.Range("H4").Value = 0.01956516564143236
$ws.Range("I4").Value = 0.867693475396689
$ws.Range("J4").Value = 0.01673637578672741
$ws.Range("K4").Value = 0.8727734419183349
$ws.Range("L4").Value = 0.01201148350833514

# Row 5: DTREE (was CART)
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.8794597600699132
$ws.Range("D5").Value = 0.01154013604599403
$ws.Range("E5").Value = 0.8747267503591282
$ws.Range("F5").Value = 0.01430326181444527
$ws.Range("G5").Value = 0.8727869765960834
$ws.Range("H5").Value = 0.01465708651317786
$ws.Range("I5").Value = 0.8773450478318153
$ws.Range("J5").Value = 0.01116752617209989
$ws.Range("K5").Value = 0.8722944237077254
$ws.Range("L5").Value = 0.007756545282219045

# Row 6: RTREE
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8729554302057678
$ws.Range("D6").Value = 0.01615956317998972
$ws.Range("E6").Value = 0.8711825960561205
$ws.Range("F6").Value = 0.01602052000354388
$ws.Range("G6").Value = 0.85651629456358
$ws.Range("H6").Value = 0.01105774845844906
$ws.Range("I6").Value = 0.8576588188119378
$ws.Range("J6").Value = 0.006769903707977678
$ws.Range("K6").Value = 0.8668397240131208
$ws.Range("L6").Value = 0.02219782299725197

# Row 7: XTREE
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.9111808479648313
$ws.Range("D7").Value = 0.01001850691035889
$ws.Range("E7").Value = 0.9044012355142202
$ws.Range("F7").Value = 0.01081958017106631
$ws.Range("G7").Value = 0.8997373029772329
$ws.Range("H7").Value = 0.01399342815343371
$ws.Range("I7").Value = 0.8917314687517901
$ws.Range("J7").Value = 0.01658942683965783
$ws.Range("K7").Value = 0.8941805225653205
$ws.Range("L7").Value = 0.01180130654359762

# Row 8: SVM (was NB previously at this row before the row9 SVM was deleted)
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.9103651915997986
$ws.Range("D8").Value = 0.006369422155240269
$ws.Range("E8").Value = 0.9027151787691416
$ws.Range("F8").Value = 0.008289604230862552
$ws.Range("G8").Value = 0.8986895766230265
$ws.Range("H8").Value = 0.01381416869271208
$ws.Range("I8").Value = 0.8955146932462622
$ws.Range("J8").Value = 0.01551327766402548
$ws.Range("K8").Value = 0.8939384684990385
$ws.Range("L8").Value = 0.006441704031766589
